$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.462.03"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "3.463.55"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'571.07"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").Value = "'160.19"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "3.465.94"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'0.574"
$ws.Range("E9").Value = "  -5.48%  "
$ws.Range("D10").Value = "'7.23"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("D12").Value = "'0.435"
$ws.Range("E12").Value = "  -3.12%  "
$ws.Range("D13").Value = "4.069.54"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D15").Value = "'27.52"
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("D16").Value = "'0.0000177"
$ws.Range("E16").Value = "  -7.44%  "
$ws.Range("D17").Value = "64.694.17"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "3.457.31"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "'6.20"
$ws.Range("E19").Value = "  -4.04%  "
$ws.Range("D20").Value = "'13.81"
$ws.Range("E20").Value = "  -3.31%  "
$ws.Range("D21").Value = "'380.19"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").Value = "'7.95"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'72.69"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("E25").Value = "  -4.96%  "
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").Value = "'6.12"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("E31").Value = "  -5.96%  "
$ws.Range("D32").Value = "'2.00"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("D33").Value = "'23.30"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("D34").Value = "'7.02"
$ws.Range("E34").Value = "  -3.43%  "
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").Value = "'161.29"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -3.32%  "
$ws.Range("E38").Value = "  +5.85%  "
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "'0.0744"
$ws.Range("E40").Value = "  -4.97%  "
$ws.Range("D41").Value = "2.832.50"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("E42").Value = "  -4.37%  "
$ws.Range("D43").Value = "'42.81"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("E44").Value = "  -6.33%  "
$ws.Range("D45").Value = "'25.79"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("D46").Value = "'0.0308"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("E47").Value = "  +8.20%  "
$ws.Range("E48").Value = "  +3.39%  "
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("E51").Value = "  -3.85%  "
